$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: the 9 SmartScore cells were stored as text ("0.592") and are
#     now numeric values (0.592) ---
$ws.Range("G10").Value = 0.592
$ws.Range("J10").Value = 0.574
$ws.Range("M10").Value = 0.506
$ws.Range("P10").Value = 0.636
$ws.Range("S10").Value = 0.581
$ws.Range("V10").Value = 0.549
$ws.Range("Y10").Value = 0.703
$ws.Range("AB10").Value = 0.699
$ws.Range("AE10").Value = 0.675

# --- Row 11: new submission from Streamlit app (user "jOSE") ---
$ws.Range("A11").Value = "jOSE"
$ws.Range("B11").Value = 23
$ws.Range("C11").Value = "Male"
$ws.Range("D11").Value = "2025-10-31 00:47:48"
$ws.Range("E11").Value = "{`n  ""portion"": 0.8,`n  ""diet"": 0.2857142857142857,`n  ""salt"": 0.8,`n  ""fat"": 0.2,`n  ""natural"": 0.6,`n  ""convenience"": 0.4,`n  ""price"": 0.2`n}"

# Instant Noodles · Top 1
$ws.Range("F11").Value = "Nongshim Neoguri Spicy Seafood"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "0.588"
$ws.Range("H11").Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"

# Instant Noodles · Top 2
$ws.Range("I11").Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Range("J11").NumberFormat = "@"
$ws.Range("J11").Value = "0.538"
$ws.Range("K11").Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"

# Instant Noodles · Top 3
$ws.Range("L11").Value = "Maruchan Ramen Sabor Pollo"
$ws.Range("M11").NumberFormat = "@"
$ws.Range("M11").Value = "0.433"
$ws.Range("N11").Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"

# Mac & Cheese · Top 1
$ws.Range("O11").Value = "Amy’s Macaroni & Cheese (frozen)"
$ws.Range("P11").NumberFormat = "@"
$ws.Range("P11").Value = "0.729"
$ws.Range("Q11").Value = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"

# Mac & Cheese · Top 2
$ws.Range("R11").Value = "Kraft Macaroni & Cheese Dinner"
$ws.Range("S11").NumberFormat = "@"
$ws.Range("S11").Value = "0.630"
$ws.Range("T11").Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"

# Mac & Cheese · Top 3
$ws.Range("U11").Value = "Annie’s Shells & White Cheddar"
$ws.Range("V11").NumberFormat = "@"
$ws.Range("V11").Value = "0.608"
$ws.Range("W11").Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"

# Ready to Eat · Top 1
$ws.Range("X11").Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Range("Y11").NumberFormat = "@"
$ws.Range("Y11").Value = "0.733"
$ws.Range("Z11").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"

# Ready to Eat · Top 2
$ws.Range("AA11").Value = "Kitchens of India Variety Pack"
$ws.Range("AB11").NumberFormat = "@"
$ws.Range("AB11").Value = "0.599"
$ws.Range("AC11").Value = "Sabor auténtico, variedad, vegetariano, necesita arroz o pan, buena calidad"

# Ready to Eat · Top 3
$ws.Range("AD11").Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Range("AE11").NumberFormat = "@"
$ws.Range("AE11").Value = "0.518"
$ws.Range("AF11").Value = "Portátil, saludable, fácil, buena textura, sabor suave"
